# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # row -> F value, G value (G omitted when unchanged)
    $ws.Range("F2").Value = 6610
    $ws.Range("G2").Value = 65

    if ($sheetName -eq "展览") {
        $ws.Range("F12").Value = 391
        $ws.Range("F13").Value = 1277
        $ws.Range("F15").Value = 3287
        $ws.Range("F19").Value = 39
        $ws.Range("F21").Value = 123
    }
    else {
        $ws.Range("F13").Value = 391
        $ws.Range("F14").Value = 1277
        $ws.Range("F16").Value = 3287
        $ws.Range("F20").Value = 39
        $ws.Range("F22").Value = 123
    }
}
